# Update the player roster table (A2:C19) on Sheet1 to reflect the new
# player list, positions and teams.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jordan Poole",             "PG,SG",    "Washington Wizards"),
    @("Jamal Murray",             "PG,SG",    "Denver Nuggets"),
    @("Shai Gilgeous-Alexander",  "PG,SG",    "Oklahoma City Thunder"),
    @("Jordan Clarkson",          "SG,SF",    "Utah Jazz"),
    @("Grayson Allen",            "PG,SG,SF", "Phoenix Suns"),
    @("Toumani Camara",           "SF,PF",    "Portland Trail Blazers"),
    @("Lauri Markkanen",          "SF,PF",    "Utah Jazz"),
    @("John Collins",             "PF,C",     "Utah Jazz"),
    @("Bam Adebayo",              "C",        "Miami Heat"),
    @("Wendell Carter Jr.",       "PF,C",     "Orlando Magic"),
    @("Nikola Jovic",             "PF,C",     "Miami Heat"),
    @("CJ McCollum",              "PG,SG",    "New Orleans Pelicans"),
    @("Tobias Harris",            "SF,PF",    "Detroit Pistons"),
    @("Kyrie Irving",             "PG,SG",    "Dallas Mavericks"),
    @("Jose Alvarado",            "PG",       "New Orleans Pelicans"),
    @("Zach LaVine",              "SG,SF",    "Chicago Bulls"),
    @("Joel Embiid",              "C",        "Philadelphia 76ers"),
    @("Mike Conley",              "PG",       "Minnesota Timberwolves")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}
